# Weekly update: a new price record (week of 2021-09-21, serial 44460) is
# inserted at the top of the "Vega Monumental Concepción - Coliflor" data
# block (rows 101-102), pushing every existing record down by two rows.
# The sheet's used range therefore grows from A1:R141 to A1:R143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 101. Excel shifts
# all rows 101:141 down to 103:143 (formatting/styles carried over from the
# row above, matching how this workbook already looks for every other
# weekly insert).
$ws.Rows("101:102").Insert()

# New "Primera" quality record
$ws.Range("A101").Value = 11
$ws.Range("B101").Value = "Vega Monumental Concepción"
$ws.Range("C101").Value = "Bíobío"
$ws.Range("D101").Value = 44460
$ws.Range("E101").Value = 8
$ws.Range("F101").Value = 100112008
$ws.Range("G101").Value = "Coliflor"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 1000
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 700
$ws.Range("M101").Value = 650
$ws.Range("N101").Value = "`$/unidad"
$ws.Range("O101").Value = "Región Metropolitana"
$ws.Range("P101").Value = 650
$ws.Range("Q101").Value = 1
$ws.Range("R101").Value = "Hortaliza"

# New "Segunda" quality record (same week)
$ws.Range("A102").Value = 11
$ws.Range("B102").Value = "Vega Monumental Concepción"
$ws.Range("C102").Value = "Bíobío"
$ws.Range("D102").Value = 44460
$ws.Range("E102").Value = 8
$ws.Range("F102").Value = 100112008
$ws.Range("G102").Value = "Coliflor"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Segunda"
$ws.Range("J102").Value = 500
$ws.Range("K102").Value = 500
$ws.Range("L102").Value = 500
$ws.Range("M102").Value = 500
$ws.Range("N102").Value = "`$/unidad"
$ws.Range("O102").Value = "Región Metropolitana"
$ws.Range("P102").Value = 500
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = "Hortaliza"
